# Insert a new weekly price record for Cilantro (Femacal de La Calera) as
# row 344, pushing the existing rows 344:397 down to 345:398.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 344 (shifts rows 344-397 -> 345-398,
# copying formatting - including the date style on column D - from the
# row above, same as a manual "Insert" in Excel).
$ws.Rows.Item(344).Insert()

# Populate the newly inserted row with the new data record.
$ws.Cells.Item(344, 1).Value  = 3
$ws.Cells.Item(344, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(344, 3).Value  = "Coquimbo"
$ws.Cells.Item(344, 4).Value  = 44776
$ws.Cells.Item(344, 5).Value  = 5
$ws.Cells.Item(344, 6).Value  = 100112040
$ws.Cells.Item(344, 7).Value  = "Cilantro"
$ws.Cells.Item(344, 8).Value  = "Sin especificar"
$ws.Cells.Item(344, 9).Value  = "Primera"
$ws.Cells.Item(344, 10).Value = 105
$ws.Cells.Item(344, 11).Value = 4500
$ws.Cells.Item(344, 12).Value = 5000
$ws.Cells.Item(344, 13).Value = 4738
$ws.Cells.Item(344, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(344, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(344, 16).Value = 1579
$ws.Cells.Item(344, 17).Value = 3
$ws.Cells.Item(344, 18).Value = "Hortaliza"
